$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New "Opportunities" sheet, placed right after "Campaigns"
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Opportunities"

function Format-Header($range) {
    $range.Font.Bold = $true
    $range.Interior.Color = 65535
    $range.Borders.LineStyle = 1
    $range.VerticalAlignment = -4108
}
function Format-Data($range) {
    $range.Borders.LineStyle = 1
    $range.VerticalAlignment = -4108
}
function Format-DataWrap($range) {
    $range.Borders.LineStyle = 1
    $range.VerticalAlignment = -4108
    $range.WrapText = $true
}

# ---- Block 1 : rows 1-2 ----
$ws2.Cells.Item(1,1).Value = "TC_ID"
$ws2.Cells.Item(1,2).Value = "TestcaseName"
$ws2.Cells.Item(1,3).Value = "OpportunityName"
$ws2.Cells.Item(1,4).Value = "RelatedTo"
$ws2.Cells.Item(1,5).Value = "ContactName"
Format-Header($ws2.Range("A1:E1"))

$ws2.Cells.Item(2,1).Value = "TC_002"
$ws2.Cells.Item(2,2).Value = "Create_opportunity_with_Contact"
$ws2.Cells.Item(2,3).Value = "Client1"
$ws2.Cells.Item(2,4).Value = "Contacts"
$ws2.Cells.Item(2,5).Value = "Asha89"
Format-Data($ws2.Range("A2:E2"))

# ---- Block 2 : rows 4-5 ----
$ws2.Cells.Item(4,1).Value = "TC_ID"
$ws2.Cells.Item(4,2).Value = "TestcaseName"
$ws2.Cells.Item(4,3).Value = "OpportunityName"
$ws2.Cells.Item(4,4).Value = "RelatedTo"
$ws2.Cells.Item(4,5).Value = "ContactName"
$ws2.Cells.Item(4,6).Value = "Subject"
$ws2.Cells.Item(4,7).Value = "OrganizationName"
$ws2.Cells.Item(4,8).Value = "BillingAddress"
$ws2.Cells.Item(4,9).Value = "ProductName"
$ws2.Cells.Item(4,10).Value = "Qty"
Format-Header($ws2.Range("A4:J4"))

$ws2.Cells.Item(5,1).Value = "TC_012"
$ws2.Cells.Item(5,2).Value = "Create_Invoice_For_Opportunity"
$ws2.Cells.Item(5,3).Value = "Client1"
$ws2.Cells.Item(5,4).Value = "Contacts"
$ws2.Cells.Item(5,5).Value = "Asha89"
$ws2.Cells.Item(5,6).Value = "Demo"
$ws2.Cells.Item(5,7).Value = "Instagram_31"
$ws2.Cells.Item(5,8).Value = "3rd floor Gopalan coworks, `nKathriguppe"
$ws2.Cells.Item(5,9).Value = "Volvo"
$ws2.Cells.Item(5,10).Value = 1
Format-Data($ws2.Range("A5:G5"))
Format-DataWrap($ws2.Range("H5"))
Format-Data($ws2.Range("I5:J5"))
$ws2.Rows.Item(5).RowHeight = 55.5

# ---- Block 3 : rows 7-8 ----
$ws2.Cells.Item(7,1).Value = "TC_ID"
$ws2.Cells.Item(7,2).Value = "TestcaseName"
$ws2.Cells.Item(7,3).Value = "OpportunityName"
$ws2.Cells.Item(7,4).Value = "RelatedTo"
$ws2.Cells.Item(7,5).Value = "ContactName"
$ws2.Cells.Item(7,6).Value = "Title"
Format-Header($ws2.Range("A7:F7"))

$ws2.Cells.Item(8,1).Value = "TC_013"
$ws2.Cells.Item(8,2).Value = "Add_Document_to_Opportunity"
$ws2.Cells.Item(8,3).Value = "Client1"
$ws2.Cells.Item(8,4).Value = "Contacts"
$ws2.Cells.Item(8,5).Value = "Asha89"
$ws2.Cells.Item(8,6).Value = "Title2"
Format-Data($ws2.Range("A8:F8"))

$ws2.PageSetup.Orientation = 1

# Column widths (bestFit in the source file)
$ws2.Columns.Item(2).ColumnWidth = 31.85546875
$ws2.Columns.Item(3).ColumnWidth = 17.42578125
$ws2.Columns.Item(4).ColumnWidth = 10
$ws2.Columns.Item(5).ColumnWidth = 13.140625
$ws2.Columns.Item(6).ColumnWidth = 7.5703125
$ws2.Columns.Item(7).ColumnWidth = 17.85546875
$ws2.Columns.Item(8).ColumnWidth = 23.42578125
$ws2.Columns.Item(9).ColumnWidth = 13.28515625

[void]$ws2.Activate()
[void]$ws2.Range("G11").Select()
